$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: update the zh-cn "Latest HO Xliff Generate Date" for the first file row.
$wsOverview.Range("G2").Value = "2016-08-18 10:51:32"

# zh-cn sheet: update Correspond Handoff Datetime and Correspond Handback DateTime for the first file row.
$wsZhCn.Range("H2").Value = "2016-08-18 10:51:27"
$wsZhCn.Range("K2").Value = "2016-08-18 10:51:42"

# de-de sheet: update Correspond Handoff Datetime and Correspond Handback DateTime for the first file row.
$wsDeDe.Range("H2").Value = "2016-08-18 10:51:32"
$wsDeDe.Range("K2").Value = "2016-08-18 10:51:50"
